$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update libraryPreparer (column B) and purpose (column E) for all data rows.
$ws.Range("B2:B35").Value = "H.BROWN"
$ws.Range("E2:E35").Value = "fullRNASEQ"

# Match the author's saved selection state.
$ws.Range("E2:E12").Select()
